$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.140.81'
$ws.Range("E2").Value = '  -1.84%  '
$ws.Range("D3").Value = '2.442.61'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''580.13'
$ws.Range("E5").Value = '  +1.21%  '
$ws.Range("D6").Value = '''142.96'
$ws.Range("E6").Value = '  -3.22%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '''0.529'
$ws.Range("E8").Value = '  -2.14%  '
$ws.Range("D9").Value = '2.444.66'
$ws.Range("E9").Value = '  -0.42%  '
$ws.Range("E10").Value = '  -3.98%  '
$ws.Range("E11").Value = '  +1.33%  '
$ws.Range("D12").Value = '''5.19'
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("D13").Value = '''0.345'
$ws.Range("E13").Value = '  -2.88%  '
$ws.Range("D14").Value = '''26.35'
$ws.Range("E14").Value = '  -3.30%  '
$ws.Range("D15").Value = '''0.0000172'
$ws.Range("E15").Value = '  -4.71%  '
$ws.Range("D16").Value = '2.864.97'
$ws.Range("E16").Value = '  -1.45%  '
$ws.Range("D17").Value = '62.102.04'
$ws.Range("E17").Value = '  -1.98%  '
$ws.Range("D18").Value = '2.436.12'
$ws.Range("E18").Value = '  -1.13%  '
$ws.Range("D19").Value = '''10.92'
$ws.Range("E19").Value = '  -3.96%  '
$ws.Range("D20").Value = '''7.09'
$ws.Range("E20").Value = '  -4.07%  '
$ws.Range("D21").Value = '''329.13'
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("D22").Value = '''4.11'
$ws.Range("E22").Value = '  -2.26%  '
$ws.Range("D23").Value = '''1.95'
$ws.Range("E23").Value = '  -5.85%  '
$ws.Range("E24").Value = '  -4.05%  '
$ws.Range("D25").Value = '''65.64'
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("D26").Value = '''9.39'
$ws.Range("E26").Value = '  +4.26%  '
$ws.Range("D27").Value = '''620.19'
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("D28").Value = '2.560.37'
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0946'
$ws.Range("E30").Value = '  -8.94%  '
$ws.Range("D31").Value = '''1.42'
$ws.Range("E31").Value = '  -6.22%  '
$ws.Range("D32").Value = '''7.99'
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("D33").Value = '''0.143'
$ws.Range("E33").Value = '  +1.14%  '
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").Value = '''4.91'
$ws.Range("E35").Value = '  -5.85%  '
$ws.Range("E36").Value = '  +0.14%  '
$ws.Range("D37").Value = '''1.43'
$ws.Range("E37").Value = '  -6.49%  '
$ws.Range("D38").Value = '''0.375'
$ws.Range("E38").Value = '  -1.65%  '
$ws.Range("D39").Value = '''149.99'
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("D40").Value = '''18.30'
$ws.Range("E40").Value = '  -2.74%  '
$ws.Range("D41").Value = '''5.24'
$ws.Range("E41").Value = '  -4.45%  '
$ws.Range("D42").Value = '''1.76'
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("D43").Value = '''42.80'
$ws.Range("E43").Value = '  +2.26%  '
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '''2.47'
$ws.Range("E45").Value = '  -9.20%  '
$ws.Range("D46").Value = '''142.83'
$ws.Range("E46").Value = '  -4.57%  '
$ws.Range("D47").Value = '''3.63'
$ws.Range("E47").Value = '  -3.89%  '
$ws.Range("D48").Value = '''0.0524'
$ws.Range("E48").Value = '  -2.87%  '
$ws.Range("D49").Value = '''0.599'
$ws.Range("E49").Value = '  -0.57%  '
$ws.Range("D50").Value = '''19.55'
$ws.Range("E50").Value = '  -8.55%  '
$ws.Range("D51").Value = '0.0₆0234'
$ws.Range("E51").Value = '  +4.00%  '
